# Add data for 2022-08-18 (carjacking by month YoY)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name to reflect the new "through" date
$ws.Name = "Through 2022-08-10"

# Update the August row label to reflect the new "through" date
$ws.Range("A9").Value = "August (through 08-10)"

# Update August row (row 9) values
$ws.Range("C9").Value = 19
$ws.Range("D9").Value = 24
$ws.Range("G9").Value = 65
$ws.Range("H9").Value = 69
$ws.Range("I9").Value = 52

# Update Total row (row 10) values
$ws.Range("C10").Value = 321
$ws.Range("D10").Value = 489
$ws.Range("G10").Value = 686
$ws.Range("H10").Value = 979
$ws.Range("I10").Value = 1022
